# Update the "想去人数" (F column) counts on the 展览 (Exhibition) and
# 全部类型 (All types) sheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# 展览 sheet updates: row -> new F value
$exhibitionUpdates = @{
    5  = 246
    11 = 1421
    12 = 37388
    13 = 7579
    14 = 129
    18 = 34
    24 = 490
    28 = 183
    30 = 438
    35 = 767
    38 = 146
    39 = 822
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# 全部类型 sheet updates: row -> new F value
$allTypesUpdates = @{
    6  = 246
    11 = 1421
    17 = 7579
    30 = 490
    33 = 183
    35 = 438
    40 = 767
    44 = 146
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}

$wb.Save()
